$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.105.08"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.444.55"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'581.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'142.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "2.440.94"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "'5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "'26.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "2.873.47"
$ws.Range("D17").Value = "62.106.12"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "2.435.72"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").Value = "'7.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'325.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "'4.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D24").Value = "'1.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").Value = "'65.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'9.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'597.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.73%  "
$ws.Range("D28").Value = "0.0₃0965"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "2.565.00"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "'7.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'4.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "'0.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'152.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("D40").Value = "'18.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'5.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "'43.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "'1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'2.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0274"
$ws.Range("E46").Value = "  +22.31%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'141.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "'3.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("D49").Value = "'0.600"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "'0.0517"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'19.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
